# Updates the "cryptos" price/volume table with freshly scraped values.
#
# Note: every data cell in the sheet is stored as plain text (the source
# workbook uses inlineStr for everything, including values that look like
# numbers, e.g. "68.373.11", "1.00", "591.62"). Assigning a numeric-looking
# string straight to Range.Value makes Excel coerce it into a real number,
# which both changes the cell type and mangles the text (loses trailing
# zeros, introduces binary floating point artifacts, etc). To keep these
# "Price" column cells as text (matching the original file byte-for-byte)
# we prefix the value with a leading apostrophe - Excel's classic "force
# text" marker - and then reset the cell style back to "Normal" so no
# stray numbering-format/quote-prefix style sticks to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.326.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "'2.512.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'591.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").Value = "'176.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "'2.511.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  +3.36%  "
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "'4.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "'0.337"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "'2.988.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "'25.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "'68.316.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "'2.494.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "'11.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").Value = "'7.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("D21").Value = "'350.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  +3.62%  "
$ws.Range("D23").Value = "'71.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.19%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "'9.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("B27").Value = "SuiNetwork"
$ws.Range("C27").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D27").Value = "'1.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.78%  "
$ws.Range("D28").Value = "'2.637.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "'0.0₃0899"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("D31").Value = "'510.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").Value = "'7.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("D37").Value = "'161.69"
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "'18.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.19%  "
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.39%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("D46").Value = "'151.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.22%  "
$ws.Range("E47").Value = "  +2.17%  "
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("D51").Value = "'0.577"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.22%  "
